# penambahan table log dan spec
# - perubahan selection/active sheet di "saldo"
# - penambahan sheet baru "Log" (table LOG001 / Log Deposit) sebagai sheet aktif

$wb = $excel.ActiveWorkbook

# --- saldo: sheet ini sebelumnya adalah tab aktif; lepaskan seleksi lamanya ---
$saldo = $wb.Worksheets.Item("saldo")
$saldo.Activate()
$saldo.Range("A1:B1").Select() | Out-Null

# --- tambahkan sheet "Log" baru setelah sheet terakhir ("saldo") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$logSheet = $wb.Worksheets.Add($null, $lastSheet)
$logSheet.Name = "Log"

# isi header + baris data tabel baru (table LOG001 / Log Deposit)
$logSheet.Range("A1").Value = "Nama Tabel"
$logSheet.Range("B1").Value = "Kode Tabel"
$logSheet.Range("A2").Value = "LOG001"
$logSheet.Range("B2").Value = "Log Deposit"

# sheet "Log" menjadi sheet aktif dengan selection di A3
$logSheet.Activate()
$logSheet.Range("A3").Select() | Out-Null
